# Update the cached "datetimeFigureOut" field text on the slide master and
# every slide layout (the date placeholder), from 28.01.2021 to 01.02.2021.
# Identify the placeholder robustly via PlaceholderFormat.Type (16 = date),
# since the current cached text can't always be trusted for a read-compare.
$p = $ppt.ActivePresentation

$newDate = "01.02.2021"
$ppPlaceholderDate = 16

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide 6 ("Wie bekommen wir Daten zu unseren Autos?") content placeholder:
# tidy up the Actor "PoissonClock" bullet text.
$slide = $p.Slides.Item(6)
$content = $slide.Shapes.Item(4)
$tf = $content.TextFrame

# Try to mark a (subtle) shrink-on-overflow line-spacing reduction, matching
# the author's autofit tweak, on a best-effort basis.
try { $content.TextFrame2.FontScale = 90 } catch {}
try { $content.TextFrame2.LineSpaceReduction = 10 } catch {}

$tr = $tf.TextRange

$full = $tr.Text
$idx = $full.IndexOf("PoissonClock")
if ($idx -ge 0) {
    $tr.Characters($idx + 1, "PoissonClock".Length).Text = "DiscreteClock"
}

$full = $tr.Text
$old2 = "random Zeiten "
$idx2 = $full.IndexOf($old2)
if ($idx2 -ge 0) {
    $tr.Characters($idx2 + 1, $old2.Length).Text = "einer bestimmten Zeiten "
}
